# Tamplo5 automation test data — "comit last 23 1 2019"
#   - rename the "login name" header to "username"
#   - drop the extra nitin.tajane@yandex.com login row (row 4)
#   - the removed header's B1 cell no longer carries the numeric/text style
#   - leave the active selection on B3 (last surviving data row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header: "login name" -> "username"
$ws.Range("A1").Value = "username"

# B1 ("password" header) loses the custom style it inherited from column B
$ws.Range("B1").Style = "Normal"

# Remove the 4th login row (nitin.tajane@yandex.com / 1234)
$ws.Rows.Item(4).Delete()

# Match the saved selection from the edited workbook
$null = $ws.Range("B3").Select()
